# Applies the "Updated cryptos list" data refresh to the Price (D) and
# Volume(1h) (E) columns, plus the WrappedBTC/Avalanche row swap (rows 20-21,
# columns B/C), exactly as captured in the source OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.889.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.843.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.22"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4748"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3668"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07193"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9268"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.63"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07674"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.897.59"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.305"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.393"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.64"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008624"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.56"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.98%  "
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.922.87"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.63"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.916"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.19"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.13"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.000"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.20"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.936"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08856"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.292"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7501"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.170"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.481"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.708"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01950"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05257"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.958"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5200"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.960"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.03%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.210"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.39%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4723"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.007"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.61"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.600"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.21"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06022"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8869"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.51%  "
